$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hyperlinks clean-up -------------------------------------------------
# Remove only the stale D3 hyperlink (old emp20@mycompany.com address).
# E2 / E3 hyperlinks (the "Cc" column) are left untouched so their
# relationship ids / uids stay intact.
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address(0, 0)
    if ($addr -eq "D3") {
        $hl.Delete()
    }
}

# --- Row 2: Investor 1 / Emp1 --------------------------------------------
$ws.Range("C2").Value = "I1"
$ws.Range("D2").Value = "emp1@investor1.com"

# --- Row 3: Investor 1 / Emp2 --------------------------------------------
$ws.Range("C3").Value = "I1"
$ws.Range("D3").Value = "emp2@investor1.com"

# --- Row 4 (new): Investor 2 / Emp1 --------------------------------------
$ws.Range("A4").Value = "Investor 2"
$ws.Range("B4").Value = "Emp1"
$ws.Range("C4").Value = "I2"
$ws.Range("D4").Value = "emp1@investor2.com"
$ws.Range("E4").Value = "advisor1@gmail.com,advisor2@gmail.com"
$ws.Range("F4").Value = "IN(91)"
$ws.Range("G4").Value = 9999999999
$ws.Range("H4").Value = "Yes"
$ws.Range("I4").Value = "Yes"
$ws.Range("J4").Value = "No"

# --- New hyperlinks for the Email column (D) and the new Cc cell (E4) ----
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:emp1@investor1.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:emp2@investor1.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:emp1@investor2.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:advisor1@gmail.com,advisor2@gmail.com") | Out-Null

# The Cc column never carries the visual "Hyperlink" cell style (E2/E3
# don't have it either), so put E4 back to the plain/default look.
$ws.Range("E4").Style = "Normal"

# --- Selection shown when the workbook is reopened ------------------------
$ws.Range("D5").Select()
